$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 12.40257566666667
$ws.Range("H2").Value = 37.20772700000001
$ws.Range("I2").Value = 0.3428195387305676
$ws.Range("J2").Value = 0.3428195387305676
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.217721
$ws.Range("N2").Value = 9.653162999999999
$ws.Range("O2").Value = 0.2700934312193076
$ws.Range("P2").Value = 0.2700934312193076
$ws.Range("Q2").Value = 39.90802817672234
$ws.Range("R2").Value = 359.172253590501
$ws.Range("S2").Value = 0.0925933055047593
$ws.Range("T2").Value = 0.09259330550475932

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 12.40257566666667
$ws.Range("H3").Value = 37.20772700000001
$ws.Range("I3").Value = 0.3428195387305676
$ws.Range("J3").Value = 0.3428195387305676
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.695641666666667
$ws.Range("N3").Value = 26.086925
$ws.Range("O3").Value = 0.7299065687806925
$ws.Range("P3").Value = 0.7299065687806925
$ws.Range("Q3").Value = 107.8483537410528
$ws.Range("R3").Value = 970.6351836694752
$ws.Range("S3").Value = 0.2502262332258083
$ws.Range("T3").Value = 0.2502262332258083

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.65713966666667
$ws.Range("H4").Value = 31.971419
$ws.Range("I4").Value = 0.294573950033059
$ws.Range("J4").Value = 0.2945739500330591
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.217721
$ws.Range("N4").Value = 9.653162999999999
$ws.Range("O4").Value = 0.2700934312193076
$ws.Range("P4").Value = 0.2700934312193076
$ws.Range("Q4").Value = 34.29170210536633
$ws.Range("R4").Value = 308.625318948297
$ws.Range("S4").Value = 0.07956248891225377
$ws.Range("T4").Value = 0.07956248891225379

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 10.65713966666667
$ws.Range("H5").Value = 31.971419
$ws.Range("I5").Value = 0.294573950033059
$ws.Range("J5").Value = 0.2945739500330591
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.695641666666667
$ws.Range("N5").Value = 26.086925
$ws.Range("O5").Value = 0.7299065687806925
$ws.Range("P5").Value = 0.7299065687806925
$ws.Range("Q5").Value = 92.67066773295278
$ws.Range("R5").Value = 834.036009596575
$ws.Range("S5").Value = 0.2150114611208053
$ws.Range("T5").Value = 0.2150114611208053

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.916194666666667
$ws.Range("H6").Value = 8.748584000000001
$ws.Range("I6").Value = 0.08060652378538533
$ws.Range("J6").Value = 0.08060652378538534
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.217721
$ws.Range("N6").Value = 9.653162999999999
$ws.Range("O6").Value = 0.2700934312193076
$ws.Range("P6").Value = 0.2700934312193076
$ws.Range("Q6").Value = 9.383500819021334
$ws.Range("R6").Value = 84.45150737119201
$ws.Range("S6").Value = 0.02177129258785545
$ws.Range("T6").Value = 0.02177129258785545

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.916194666666667
$ws.Range("H7").Value = 8.748584000000001
$ws.Range("I7").Value = 0.08060652378538533
$ws.Range("J7").Value = 0.08060652378538534
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.695641666666667
$ws.Range("N7").Value = 26.086925
$ws.Range("O7").Value = 0.7299065687806925
$ws.Range("P7").Value = 0.7299065687806925
$ws.Range("Q7").Value = 25.35818385157778
$ws.Range("R7").Value = 228.2236546642
$ws.Range("S7").Value = 0.05883523119752988
$ws.Range("T7").Value = 0.05883523119752989

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.547057666666666
$ws.Range("H8").Value = 16.641173
$ws.Range("I8").Value = 0.1533261962440107
$ws.Range("J8").Value = 0.1533261962440107
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.217721
$ws.Range("N8").Value = 9.653162999999999
$ws.Range("O8").Value = 0.2700934312193076
$ws.Range("P8").Value = 0.2700934312193076
$ws.Range("Q8").Value = 17.84888394224433
$ws.Range("R8").Value = 160.639955480199
$ws.Range("S8").Value = 0.04141239843934975
$ws.Range("T8").Value = 0.04141239843934976

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.547057666666666
$ws.Range("H9").Value = 16.641173
$ws.Range("I9").Value = 0.1533261962440107
$ws.Range("J9").Value = 0.1533261962440107
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.695641666666667
$ws.Range("N9").Value = 26.086925
$ws.Range("O9").Value = 0.7299065687806925
$ws.Range("P9").Value = 0.7299065687806925
$ws.Range("Q9").Value = 48.23522577366944
$ws.Range("R9").Value = 434.117031963025
$ws.Range("S9").Value = 0.1119137978046609
$ws.Range("T9").Value = 0.1119137978046609

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.655179333333333
$ws.Range("H10").Value = 13.965538
$ws.Range("I10").Value = 0.1286737912069773
$ws.Range("J10").Value = 0.1286737912069773
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.217721
$ws.Range("N10").Value = 9.653162999999999
$ws.Range("O10").Value = 0.2700934312193076
$ws.Range("P10").Value = 0.2700934312193076
$ws.Range("Q10").Value = 14.97906829963266
$ws.Range("R10").Value = 134.811614696694
$ws.Range("S10").Value = 0.03475394577508926
$ws.Range("T10").Value = 0.03475394577508927

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4.655179333333333
$ws.Range("H11").Value = 13.965538
$ws.Range("I11").Value = 0.1286737912069773
$ws.Range("J11").Value = 0.1286737912069773
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 8.695641666666667
$ws.Range("N11").Value = 26.086925
$ws.Range("O11").Value = 0.7299065687806925
$ws.Range("P11").Value = 0.7299065687806925
$ws.Range("Q11").Value = 40.47977137673889
$ws.Range("R11").Value = 364.31794239065
$ws.Range("S11").Value = 0.09391984543188804
$ws.Range("T11").Value = 0.09391984543188807
